$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This report tracks two source files:
#   23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md
#   2b8b8f24-b4f7-4a6d-baad-590544370594.md
#
# A new handoff was generated for 23b6880d-...md, so its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", its handoff
# datetime is refreshed, and (to mirror the generated report's row order)
# it now appears in the row below 2b8b8f24-...md.
# ---------------------------------------------------------------------------

function Set-SheetHyperlinkText {
    param($ws, $addressToText)

    $links = @()
    foreach ($hl in $ws.Hyperlinks) {
        $links += $hl
    }
    foreach ($hl in $links) {
        $addr = $hl.Range.Address()
        if ($addressToText.ContainsKey($addr)) {
            $hl.TextToDisplay = $addressToText[$addr]
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

$wsOverview.Range("A3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

Set-SheetHyperlinkText $wsOverview @{
    "`$A`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
    "`$A`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-02-29 03:58:46"
$wsZh.Range("E2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
$wsZh.Range("F2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-02-29 03:59:48"
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("A3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-02-29 04:01:06"
$wsZh.Range("E3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
$wsZh.Range("F3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-02-29 03:59:48"
$wsZh.Range("H3").Value = "Include"

Set-SheetHyperlinkText $wsZh @{
    "`$A`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
    "`$C`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.zh-cn.xlf"
    "`$E`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
    "`$F`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.zh-cn.xlf"
    "`$A`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
    "`$C`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.zh-cn.xlf"
    "`$E`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
    "`$F`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.zh-cn.xlf"
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.de-de.xlf"
$wsDe.Range("D2").Value = "2016-02-29 03:58:57"
$wsDe.Range("E2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
$wsDe.Range("F2").Value = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.de-de.xlf"
$wsDe.Range("G2").Value = "2016-02-29 04:00:15"
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("A3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.de-de.xlf"
$wsDe.Range("D3").Value = "2016-02-29 04:01:21"
$wsDe.Range("E3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
$wsDe.Range("F3").Value = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.de-de.xlf"
$wsDe.Range("G3").Value = "2016-02-29 04:00:15"
$wsDe.Range("H3").Value = "Include"

Set-SheetHyperlinkText $wsDe @{
    "`$A`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
    "`$C`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.de-de.xlf"
    "`$E`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.md"
    "`$F`$2" = "2b8b8f24-b4f7-4a6d-baad-590544370594.70046711852f25795be2495f9edad3f1f83da7c0.de-de.xlf"
    "`$A`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
    "`$C`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.de-de.xlf"
    "`$E`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.md"
    "`$F`$3" = "23b6880d-c0b5-437c-bb68-1fd3ff3472c3.3562a1d5916b576694b1713294ffc4bfb7de1c2b.de-de.xlf"
}

Write-Host "Done"
